$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list on Thu Mar 30 20:39:36 UTC 2023 with GitHub Actions

$ws.Range("D2").Value = "28.033.32"
$ws.Range("E2").Value = "  -1.41%  "

$ws.Range("D3").Value = "1.788.90"
$ws.Range("E3").Value = "  -1.05%  "

$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'316.73"
$ws.Range("E5").Value = "  +0.09%  "

$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").Value = "'0.5349"
$ws.Range("E7").Value = "  -2.15%  "

$ws.Range("D8").Value = "'0.3766"
$ws.Range("E8").Value = "  -2.24%  "

$ws.Range("D9").Value = "'0.07455"
$ws.Range("E9").Value = "  -1.91%  "

$ws.Range("D10").Value = "'41.78"
$ws.Range("E10").Value = "  -0.98%  "

$ws.Range("D11").Value = "'1.095"
$ws.Range("E11").Value = "  -2.99%  "

$ws.Range("D12").Value = "'1.000"
$ws.Range("E12").Value = "  +0.07%  "

$ws.Range("D13").Value = "'20.61"
$ws.Range("E13").Value = "  -2.85%  "

$ws.Range("D14").Value = "'6.103"
$ws.Range("E14").Value = "  -1.54%  "

$ws.Range("D15").Value = "'7.222"
$ws.Range("E15").Value = "  -2.33%  "

$ws.Range("D16").Value = "1.776.07"
$ws.Range("E16").Value = "  -1.63%  "

$ws.Range("D17").Value = "'88.82"
$ws.Range("E17").Value = "  -3.61%  "

$ws.Range("D18").Value = "'0.00001055"
$ws.Range("E18").Value = "  -1.68%  "

$ws.Range("D19").Value = "'0.06446"
$ws.Range("E19").Value = "  -0.04%  "

$ws.Range("D20").Value = "'0.9999"
$ws.Range("E20").Value = "  +0.07%  "

$ws.Range("D21").Value = "'17.31"
$ws.Range("E21").Value = "  -0.46%  "

$ws.Range("D22").Value = "'5.899"
$ws.Range("E22").Value = "  -1.47%  "

$ws.Range("D23").Value = "28.057.88"
$ws.Range("E23").Value = "  -1.36%  "

$ws.Range("D24").Value = "'11.23"
$ws.Range("E24").Value = "  -2.29%  "

$ws.Range("E25").Value = "  -2.10%  "

$ws.Range("D26").Value = "'155.18"
$ws.Range("E26").Value = "  -2.40%  "

$ws.Range("D27").Value = "'20.28"
$ws.Range("E27").Value = "  -2.06%  "

$ws.Range("D28").Value = "1.991.30"
$ws.Range("E28").Value = "  -1.16%  "

$ws.Range("D29").Value = "'2.285"
$ws.Range("E29").Value = "  -5.09%  "

$ws.Range("D30").Value = "'120.16"
$ws.Range("E30").Value = "  -2.97%  "

$ws.Range("E31").Value = "  -1.74%  "

$ws.Range("D32").Value = "'0.1049"
$ws.Range("E32").Value = "  +2.86%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'5.555"
$ws.Range("E33").Value = "  -3.62%  "

$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'3.635"
$ws.Range("E34").Value = "  -1.46%  "

$ws.Range("E35").Value = "  -2.67%  "

$ws.Range("D36").Value = "'0.06444"
$ws.Range("E36").Value = "  -0.53%  "

$ws.Range("D37").Value = "'0.02288"
$ws.Range("E37").Value = "  -1.76%  "

$ws.Range("D38").Value = "'5.019"
$ws.Range("E38").Value = "  -2.24%  "

$ws.Range("D39").Value = "'8.470"
$ws.Range("E39").Value = "  -4.54%  "

$ws.Range("B40").Value = "WEMIXTOKEN"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").Value = "'1.447"
$ws.Range("E40").Value = "  +4.53%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.6160"
$ws.Range("E41").Value = "  -4.27%  "

$ws.Range("D42").Value = "'11.10"
$ws.Range("E42").Value = "  -4.81%  "

$ws.Range("D43").Value = "'1.170"
$ws.Range("E43").Value = "  +0.76%  "

$ws.Range("D44").Value = "'0.9998"
$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("D45").Value = "'13.25"
$ws.Range("E45").Value = "  -1.93%  "

$ws.Range("D46").Value = "'3.668"
$ws.Range("E46").Value = "  -0.47%  "

$ws.Range("D47").Value = "'0.5767"
$ws.Range("E47").Value = "  -3.72%  "

$ws.Range("D48").Value = "'127.12"
$ws.Range("E48").Value = "  +0.01%  "

$ws.Range("D49").Value = "'1.191"
$ws.Range("E49").Value = "  +3.73%  "

$ws.Range("D50").Value = "'1.926"
$ws.Range("E50").Value = "  -3.33%  "

$ws.Range("D51").Value = "'0.06809"
$ws.Range("E51").Value = "  -1.39%  "
